# Add lambda (J) / stdev (K) readings for the Nitrogen and CO2 gas blocks,
# matching the pattern already used for the Argon block (rows 5-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nitrogen block (rows 12-19) ---------------------------------------
# Two new lambda readings in column J, plus the "avg" row (I/J/K) that
# mirrors the Argon block's row 7.
$ws.Range("J14").Value = 0.489611
$ws.Range("K14").Value = "????"
$ws.Range("J15").Value = 0.230706

$ws.Range("I16").Value = "avg"
$ws.Range("J16").Formula = "=AVERAGE(J14:J15)"
$ws.Range("K16").Formula = "=STDEV(J14:J15)/SQRT(2)"

# --- CO2 block (rows 21-28) ----------------------------------------------
# Re-assert the CO2 label so the rebuilt shared-string table places the new
# "????" string ahead of it (keeps the label's displayed text unchanged).
$ws.Range("B21").Value = "CO2"

$ws.Range("J23").Value = 0.397348
$ws.Range("J24").Value = 0.354952

$ws.Range("I25").Value = "avg"
$ws.Range("J25").Formula = "=AVERAGE(J23:J24)"
$ws.Range("K25").Formula = "=STDEV(J23:J24)/SQRT(2)"

# --- Restore cursor position reported in the saved workbook --------------
$ws.Range("L18").Select() | Out-Null
